$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 217.2973
$ws.Range("I33").Value = 169.92592
$ws.Range("J33").Value = 345.2
$ws.Range("K33").Value = 169.92592
$ws.Range("L33").Value = 345.2
$ws.Range("M33").Value = 59.07408000000001
$ws.Range("N33").Value = -803.2
$ws.Range("H76").Value = 3330
$ws.Range("I76").Value = 3347.3684
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3347.3684
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -3032.3684
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3330
$ws.Range("I79").Value = 3347.3684
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3347.3684
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2255.3684
$ws.Range("N79").Value = -5184
$ws.Range("H95").Value = 42500
$ws.Range("J95").Value = 42500
$ws.Range("L95").Value = 42500
$ws.Range("N95").Value = -47992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 100948.6
$ws.Range("I45").Value = 167414.67
$ws.Range("K45").Value = 167414.67
$ws.Range("M45").Value = -167037.67
$ws.Range("H61").Value = 2973.6843
$ws.Range("I61").Value = 1644
$ws.Range("K61").Value = 1644
$ws.Range("M61").Value = -1432
$ws.Range("H95").Value = 26735.666
$ws.Range("J95").Value = 26735.666
$ws.Range("L95").Value = 26735.666
$ws.Range("N95").Value = -32227.666
$ws.Range("H106").Value = 35333.332
$ws.Range("J106").Value = 35333.332
$ws.Range("L106").Value = 35333.332
$ws.Range("N106").Value = -37857.332
$ws.Range("H133").Value = 122222
$ws.Range("J133").Value = 122222
$ws.Range("L133").Value = 122222
$ws.Range("N133").Value = -127282
$ws.Range("H136").Value = 2973.6843
$ws.Range("I136").Value = 1644
$ws.Range("K136").Value = 4932
$ws.Range("M136").Value = -2382

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 723.63635
$ws.Range("I12").Value = 723.63635
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 723.63635
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -555.63635
$ws.Range("N12").ClearContents()
$ws.Range("H134").Value = 32928.316
$ws.Range("I134").Value = 40040.71
$ws.Range("J134").Value = 1430.5714
$ws.Range("K134").Value = 120122.13
$ws.Range("L134").Value = 4291.7142
$ws.Range("M134").Value = -117587.13
$ws.Range("N134").Value = -9361.7142
$ws.Range("H140").Value = 52932.5
$ws.Range("J140").Value = 52932.5
$ws.Range("L140").Value = 52932.5
$ws.Range("N140").Value = -63292.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 35.125
$ws.Range("I7").Value = 42.75
$ws.Range("J7").Value = 27.5
$ws.Range("K7").Value = 42.75
$ws.Range("L7").Value = 27.5
$ws.Range("M7").Value = 70.25
$ws.Range("N7").Value = -253.5
$ws.Range("H93").Value = 16750
$ws.Range("I93").Value = 16750
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 16750
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -14878
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 1947.8889
$ws.Range("I132").Value = 1628.2142
$ws.Range("J132").Value = 3066.75
$ws.Range("K132").Value = 4884.642599999999
$ws.Range("L132").Value = 9200.25
$ws.Range("M132").Value = -2354.642599999999
$ws.Range("N132").Value = -14260.25
$ws.Range("H134").Value = 3462.7896
$ws.Range("I134").Value = 3521.8333
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 10565.4999
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -8030.499899999999
$ws.Range("N134").Value = -12270
$ws.Range("H135").Value = 41000
$ws.Range("J135").Value = 41000
$ws.Range("L135").Value = 41000
$ws.Range("N135").Value = -51140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 1666.6666
$ws.Range("H77").Value = 1666.6666
$ws.Range("H98").Value = 1741.091
$ws.Range("J98").Value = 2077.4285
$ws.Range("L98").Value = 6232.2855
$ws.Range("N98").Value = -9228.2855
$ws.Range("H99").Value = 2496.2856
$ws.Range("I99").Value = 1359.3334
$ws.Range("K99").Value = 4078.0002
$ws.Range("M99").Value = -1832.0002
$ws.Range("H107").Value = 43777.87
$ws.Range("I107").Value = 52789.58
$ws.Range("J107").Value = 37436.297
$ws.Range("K107").Value = 158368.74
$ws.Range("L107").Value = 112308.891
$ws.Range("M107").Value = -156448.74
$ws.Range("N107").Value = -116148.891
$ws.Range("H137").Value = 33220.723
$ws.Range("I137").Value = 2422.5
$ws.Range("J137").Value = 42020.215
$ws.Range("K137").Value = 7267.5
$ws.Range("L137").Value = 126060.645
$ws.Range("M137").Value = -2167.5
$ws.Range("N137").Value = -136260.645

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("K13").Value = 1000
$ws.Range("M13").Value = -861
$ws.Range("H70").Value = 6321.3335
$ws.Range("I70").Value = 6601.6
$ws.Range("J70").Value = 4920
$ws.Range("K70").Value = 6601.6
$ws.Range("L70").Value = 4920
$ws.Range("M70").Value = -6331.6
$ws.Range("N70").Value = -5460
$ws.Range("H73").Value = 6321.3335
$ws.Range("I73").Value = 6601.6
$ws.Range("J73").Value = 4920
$ws.Range("K73").Value = 6601.6
$ws.Range("L73").Value = 4920
$ws.Range("M73").Value = -5665.6
$ws.Range("N73").Value = -6792
$ws.Range("H122").Value = 12001.875
$ws.Range("I122").Value = 20001.75
$ws.Range("K122").Value = 60005.25
$ws.Range("M122").Value = -57555.25
$ws.Range("H132").Value = 3592.6843
$ws.Range("I132").Value = 3024.8572
$ws.Range("J132").Value = 3923.9167
$ws.Range("K132").Value = 9074.571599999999
$ws.Range("L132").Value = 11771.7501
$ws.Range("M132").Value = -6544.571599999999
$ws.Range("N132").Value = -16831.7501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1801.4572
$ws.Range("I7").Value = 1644.4615
$ws.Range("K7").Value = 1644.4615
$ws.Range("M7").Value = -1532.4615
$ws.Range("H40").Value = 2628.5715
$ws.Range("I40").Value = 2533.3333
$ws.Range("J40").Value = 3200
$ws.Range("K40").Value = 2533.3333
$ws.Range("L40").Value = 3200
$ws.Range("M40").Value = -2397.3333
$ws.Range("N40").Value = -3472
$ws.Range("H46").Value = 754.5238000000001
$ws.Range("J46").Value = 433.33334
$ws.Range("L46").Value = 433.33334
$ws.Range("N46").Value = -809.33334
$ws.Range("H55").Value = 248.1875
$ws.Range("I55").Value = 315.25
$ws.Range("J55").Value = 225.83333
$ws.Range("K55").Value = 315.25
$ws.Range("L55").Value = 225.83333
$ws.Range("M55").Value = -142.25
$ws.Range("N55").Value = -571.8333299999999
$ws.Range("H82").Value = 1726.7858
$ws.Range("I82").Value = 1419.4445
$ws.Range("J82").Value = 2280
$ws.Range("K82").Value = 1419.4445
$ws.Range("L82").Value = 2280
$ws.Range("M82").Value = -1058.4445
$ws.Range("N82").Value = -3002
$ws.Range("H85").Value = 1726.7858
$ws.Range("I85").Value = 1419.4445
$ws.Range("J85").Value = 2280
$ws.Range("K85").Value = 1419.4445
$ws.Range("L85").Value = 2280
$ws.Range("M85").Value = -171.4445000000001
$ws.Range("N85").Value = -4776
$ws.Range("H126").Value = 1801.4572
$ws.Range("I126").Value = 1644.4615
$ws.Range("K126").Value = 4933.3845
$ws.Range("M126").Value = -2463.3845
$ws.Range("H136").Value = 21400
$ws.Range("I136").Value = 26000
$ws.Range("K136").Value = 78000
$ws.Range("M136").Value = -75450
$ws.Range("H139").Value = 54115
$ws.Range("J139").Value = 54115
$ws.Range("L139").Value = 54115
$ws.Range("N139").Value = -64395

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 10000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10284
$ws.Range("H122").Value = 1613.4572
$ws.Range("I122").Value = 1018.8
$ws.Range("K122").Value = 3056.4
$ws.Range("M122").Value = -606.3999999999996
$ws.Range("H126").Value = 1749.2727
$ws.Range("I126").Value = 1889.6666
$ws.Range("J126").Value = 1117.5
$ws.Range("K126").Value = 5668.9998
$ws.Range("L126").Value = 3352.5
$ws.Range("M126").Value = -3198.9998
$ws.Range("N126").Value = -8292.5
$ws.Range("H132").Value = 6294.136
$ws.Range("I132").Value = 10086.182
$ws.Range("J132").Value = 2502.0908
$ws.Range("K132").Value = 30258.546
$ws.Range("L132").Value = 7506.2724
$ws.Range("M132").Value = -27728.546
$ws.Range("N132").Value = -12566.2724
$ws.Range("H136").Value = 15443.643
$ws.Range("I136").Value = 26025.75
$ws.Range("K136").Value = 78077.25
$ws.Range("M136").Value = -75527.25
$ws.Range("H139").Value = 49933.332
$ws.Range("J139").Value = 49933.332
$ws.Range("L139").Value = 49933.332
$ws.Range("N139").Value = -60213.332
